# Auto-generated edit script: update cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.937.44"
$ws.Range("E2").Value = "  -3.20%  "
$ws.Range("D3").Value = "'1.831.78"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D4").Value = "'0.9981"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'277.85"
$ws.Range("E5").Value = "  -7.41%  "
$ws.Range("D6").Value = "'0.9977"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.5109"
$ws.Range("E7").Value = "  -4.89%  "
$ws.Range("D8").Value = "'0.3491"
$ws.Range("E8").Value = "  -6.59%  "
$ws.Range("D9").Value = "'44.64"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").Value = "'0.06797"
$ws.Range("E10").Value = "  -4.76%  "
$ws.Range("D11").Value = "'19.88"
$ws.Range("E11").Value = "  -7.72%  "
$ws.Range("D12").Value = "'0.8066"
$ws.Range("E12").Value = "  -9.02%  "
$ws.Range("D13").Value = "'0.07796"
$ws.Range("E13").Value = "  -4.33%  "
$ws.Range("D14").Value = "'1.825.75"
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("D15").Value = "'5.077"
$ws.Range("E15").Value = "  -4.00%  "
$ws.Range("D16").Value = "'88.19"
$ws.Range("E16").Value = "  -4.48%  "
$ws.Range("D17").Value = "'0.9981"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "'14.16"
$ws.Range("E18").Value = "  -4.42%  "
$ws.Range("D19").Value = "'0.000008063"
$ws.Range("E19").Value = "  -4.98%  "
$ws.Range("D20").Value = "'0.9975"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "'25.984.20"
$ws.Range("E21").Value = "  -3.13%  "
$ws.Range("D22").Value = "'4.770"
$ws.Range("E22").Value = "  -3.91%  "
$ws.Range("D23").Value = "'10.04"
$ws.Range("E23").Value = "  -5.58%  "
$ws.Range("D24").Value = "'6.188"
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("D25").Value = "'2.367"
$ws.Range("E25").Value = "  +3.73%  "
$ws.Range("D26").Value = "'142.54"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("E27").Value = "  -4.00%  "
$ws.Range("D28").Value = "'17.20"
$ws.Range("E28").Value = "  -4.20%  "
$ws.Range("D29").Value = "'109.42"
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("D30").Value = "'4.355"
$ws.Range("E30").Value = "  -7.25%  "
$ws.Range("D31").Value = "'4.289"
$ws.Range("E31").Value = "  -7.26%  "
$ws.Range("D32").Value = "'0.08779"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").Value = "'0.04866"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").Value = "'1.164"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "'0.7281"
$ws.Range("E35").Value = "  -10.39%  "
$ws.Range("D36").Value = "'2.861"
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").Value = "'3.196"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").Value = "'0.9970"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.400"
$ws.Range("E39").Value = "  -9.98%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01847"
$ws.Range("E40").Value = "  -4.96%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5149"
$ws.Range("E41").Value = "  -15.22%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9472"
$ws.Range("E42").Value = "  -11.04%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'117.07"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.228"
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'8.044"
$ws.Range("E45").Value = "  -8.12%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'0.9968"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "'0.4510"
$ws.Range("E47").Value = "  -14.62%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1361"
$ws.Range("E48").Value = "  -8.23%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.302"
$ws.Range("E49").Value = "  -6.71%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'36.20"
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05918"
$ws.Range("E51").Value = "  -2.41%  "
